$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Name"
$ws.Range("E1").Value = "entry Time"
$ws.Range("F1").Value = "location"

# Copy header style from E1 to F1 (bold, centered, bordered header)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Update data rows - name column
$ws.Range("A2").Value = "Jaya"
$ws.Range("A3").Value = "Jaya"
$ws.Range("A4").Value = "Jaya"

# Add new location column values
$ws.Range("F2").Value = 603
$ws.Range("F3").Value = 603
$ws.Range("F4").Value = 603
